$wb = $excel.ActiveWorkbook

# --- Sheet: Kongegårdsgatan Molndal Sweden (row 5 Internal Fill Rate / Commit-Forecast,
#            row 6 & 7 Manufacturing Voluntary Turnover PY Actual / AOP) ---
$ws1 = $wb.Worksheets.Item("Kongegårdsgatan Molndal Sweden")

$ws1.Range("E5").Value = 0.5
$ws1.Range("L5").Value = 0
$ws1.Range("M5").Value = 0.5
$ws1.Range("N5").Value = 0.5
$ws1.Range("O5").Value = 0.5
$ws1.Range("P5").Value = 0.5
$ws1.Range("Q5").Value = 0.5
$ws1.Range("R5").Value = 0.5
$ws1.Range("S5").Value = 0.5
$ws1.Range("T5").Value = 0.5
$ws1.Range("U5").Value = 0.5
$ws1.Range("V5").Value = 0.5
$ws1.Range("W5").Value = 0.5

$ws1.Range("E6").Value = 0.0776
$ws1.Range("E7").Value = 0.0776

# --- Sheet: Charlotte  North Carolina (row 2 Professional Voluntary Turnover / Commit-Forecast) ---
$ws4 = $wb.Worksheets.Item("Charlotte  North Carolina")

$ws4.Range("K2").Value = 1
$ws4.Range("M2").Value = 0.0833333333333333
$ws4.Range("N2").Value = 0.25
$ws4.Range("O2").Value = 0.0833333333333333
$ws4.Range("P2").Value = 0.0833333333333333
$ws4.Range("Q2").Value = 0.0833333333333333
$ws4.Range("R2").Value = 0.25
$ws4.Range("S2").Value = 0.0833333333333333
$ws4.Range("T2").Value = 0.0833333333333333
$ws4.Range("U2").Value = 0.0833333333333333
$ws4.Range("V2").Value = 0.25
$ws4.Range("W2").Value = 1

# --- Sheet: Shanghai Minhang District China (row 3 Internal Fill Rate / Commit-Forecast) ---
$ws9 = $wb.Worksheets.Item("Shanghai Minhang District Chin")

$ws9.Range("L3").ClearContents()
